$d = $word.ActiveDocument

# --- Hunk 1: "overall summary in total of the data" -> "overall summary of the total data" ---
$d.Content.Find.Execute(
    "summary in total of the data,", $true, $false, $false, $false, $false,
    $true, 1, $false, "summary of the total data,", 2) | Out-Null

# --- Hunk 2a: "...50% of the total." -> "...50% of the total number of songs." ---
$d.Content.Find.Execute(
    "approximately 50% of the total.", $true, $false, $false, $false, $false,
    $true, 1, $false, "approximately 50% of the total number of songs.", 2) | Out-Null

# --- Hunk 2b: "is a histogram created" -> "is a chart created" ---
$d.Content.Find.Execute(
    "is a histogram created", $true, $false, $false, $false, $false,
    $true, 1, $false, "is a chart created", 2) | Out-Null

# --- Hunk 2c: "different measured. " -> "different measured features. " ---
$d.Content.Find.Execute(
    "different measured. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "different measured features. ", 2) | Out-Null

# --- Hunk 3a: "with 4,321, the following" -> "with over 4,000, the following" ---
$d.Content.Find.Execute(
    "with 4,321, the following", $true, $false, $false, $false, $false,
    $true, 1, $false, "with over 4,000, the following", 2) | Out-Null

# --- Hunk 3b: "popular songs of only 2,670" -> "popular songs with only 2,670" ---
$d.Content.Find.Execute(
    "popular songs of only 2,670", $true, $false, $false, $false, $false,
    $true, 1, $false, "popular songs with only 2,670", 2) | Out-Null

# --- Hunk 4: Instrumentalness sentence ---
$d.Content.Find.Execute(
    "our most significant, which was flat over all so that in comparison the non-popular song scores are increasing over time while the popular",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "our most significant feature, which was flat overall, but in comparison the non-popular song scores are increasing over the decades while the popular",
    2) | Out-Null

# --- Hunk 5: duration / nineties sentence ---
$d.Content.Find.Execute(
    "in the nineties the popular songs at almost 4.5 minutes were slightly",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "in the nineties the popular songs were almost 4.5 minutes and slightly",
    2) | Out-Null

# --- Hunk 6: remove the trailing paragraph (lastRenderedPageBreak + "  ") ---
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$lastPara.Range.Delete() | Out-Null

Write-Output "done"
